# Add the "Demo" worksheet (the workbook currently has no sheets at all)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Demo"

# ---- Header row ----
$ws.Range("A1").Value = "Scenario"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "Run Time"

# ---- Data rows ----
$data = @(
    @("Time taken to fully drawn (TTFD)", "24.2s", "10:44 AM EST"),
    @("State change Arizona to New York", "2.7s", "10:44 AM EST"),
    @("County Change > 1st 2 (Albany and Allegany)", "2.1s", "10:44 AM EST"),
    @("Plan type Change (Local HMO to Local PPO)", "1.6s", "10:44 AM EST"),
    @("SNP Plan Change (DSNP to Non SNP)", "2.6s", "10:44 AM EST"),
    @("Base PLAN Selection (H3418_004 to H3418_008)", "1.1s", "10:44 AM EST"),
    @("Comparison PLan Selection (H2775_105 to H2775_106)", "0.6s", "10:44 AM EST"),
    @("4 stars & up", "1.0s", "10:44 AM EST")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# ---- Header formatting: bold, centered, light-gray fill ----
# Build the combined style on A1 (font + fill + alignment in one pass so the
# engine collapses them into a single cellXf), then copy that format onto the
# rest of the header row so no intermediate/orphaned styles are left behind.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Interior.Color = 13882323
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4108

$a1.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 59.166666666666664
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668

# ---- Page setup ----
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
